$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: new "sequence to use" permutation (rows 1-20)
$seq = @(14, 12, 7, 0, 3, 1, 4, 17, 2, 16, 8, 18, 11, 6, 19, 15, 5, 10, 13, 9)
for ($i = 0; $i -lt $seq.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $seq[$i]
}

# D1: fitness improvement value
$ws.Range("D1").Value = 108.1790132698523

# D2: penalty improvement value
$ws.Range("D2").Value = 75.52659301172291

# B21: last generation fit value
$ws.Range("B21").Value = 0.8108155531546481
